$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "fonte nobreak 12v-10a fx250 metalica amfer cor preto"
$ws.Range("F2").Value = 406.6
$ws.Range("G2").Value = 406.6
$ws.Range("H2").Value = "FONTE NOBREAK FX250 12V/10A"
